# fix: Modify download template prompt
# Update the QA-pair template sheet:
#  - shorten the "segment content" column header/description text
#    (drop the "max 4096 characters" claim, since longer text can in
#    fact be imported)
#  - reset the sheet's saved cursor/selection back to B1

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MaxKB产品介绍")
$ws.Activate()

$ws.Range("B1").Value = "分段内容（必填，问题答案）"

$ws.Range("B1").Select()
